$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row => list of columns (within D..H) whose value flips from 0 to 1
$changes = @{
    3  = @("G","H")
    4  = @("D","E")
    5  = @("D","E")
    6  = @("D","E")
    7  = @("H")
    8  = @("H")
    9  = @("D","E")
    10 = @("D","E")
    11 = @("D","E")
    12 = @("D","E")
    13 = @("D","E")
    14 = @("D","E")
    15 = @("D","E")
    16 = @("H")
    17 = @("D","E")
    18 = @("D","E")
}

foreach ($row in $changes.Keys) {
    foreach ($col in $changes[$row]) {
        $ws.Range("$col$row").Value = 1
    }
}
